# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.898.69'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.505.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.519'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.504.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.51%  '
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('E12').Value = '  -4.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.11'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.958.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.58%  '
$ws.Range('E16').Value = '  -3.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.007.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.503.29'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '343.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '68.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('E29').Value = '  -3.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0977'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.22'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '522.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.88%  '
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.11'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  -3.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.80%  '
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.354'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('E42').Value = '  -3.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '146.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.553'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.43%  '
$ws.Range('E48').Value = '  -6.90%  '
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0752'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.25%  '
